$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the three new LEM/LED pairs below the existing data (rows 88-90)
$ws.Range("A88").Value = "LEM-230-33-27KH"
$ws.Range("B88").Value = "LED-230-H02-27"

$ws.Range("A89").Value = "LEM-319-00-35KH"
$ws.Range("B89").Value = "LED-319-H00-35"

$ws.Range("A90").Value = "LEM-326-00-27KS"
$ws.Range("B90").Value = "LED-326-S00-27"

# Update the selection to match the post-edit state
[void]$ws.Range("B94").Select()
